$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.143.50"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.929.67"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'591.48"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'145.14"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.506"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "'6.95"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.441"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'0.0000226"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'33.69"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "3.414.24"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "61.094.98"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "'6.72"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "2.925.37"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'437.98"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "'13.45"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "'0.678"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'7.13"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'81.74"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'10.99"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "'11.84"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").Value = "'2.60"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'7.02"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("D31").Value = "'26.68"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "'0.110"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "0.0₃0870"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'5.63"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").Value = "'1.99"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'8.61"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "'42.11"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "'0.290"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").Value = "'377.36"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'0.0347"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "2.689.92"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "'133.57"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D48").Value = "'23.95"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'2.01"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -0.06%  "
